$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: username/nama updated to new sample data
$ws.Range("B2").Value = "coba-coba"
$ws.Range("C2").Value = "admin 01"

# Row 3: username/nama updated to new sample data
$ws.Range("B3").Value = "athif"
$ws.Range("C3").Value = "athif"

# Rows 4 and 5 are no longer part of the table - remove them entirely
# so the used range shrinks back to A1:D3
$ws.Range("A4:A5").EntireRow.Delete()

# Column B now has an explicit custom width (picked up from the author's
# resize in Excel)
$ws.Columns.Item(2).ColumnWidth = 13.5

# Reflect the author's final selection/cursor position
$ws.Range("D3").Select()
